{"js": "// The document's auto-managed \"_GoBack\" bookmark currently wraps the\n// class-diagram picture (left over from whatever was typed/pasted there\n// previously). We are about to make the \"real\" last edit -- changing the\n// date from \"February 4, 2015\" to \"February 25, 2015\" -- so the bookmark\n// needs to move to mark that new edit position instead.\n\n// 1) Remove the stale \"_GoBack\" bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the date text and replace the day \"4\" with \"25\".\nconst dateHits = context.document.body.search(\"February 4, 2015\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nconst dateRange = dateHits.items[0];\nconst dayHits = dateRange.search(\"4\", { matchCase: true });\ndayHits.load(\"items\");\nawait context.sync();\n\ndayHits.items[0].insertText(\"25\", \"Replace\");\nawait context.sync();\n\n// 3) Re-locate the updated date text, and drop a collapsed \"_GoBack\"\n//    bookmark right at the point where the new text (\"25\") was inserted --\n//    i.e. immediately before the \", 2015\" that follows it.\nconst updatedHits = context.document.body.search(\"February 25, 2015\", { matchCase: true });\nupdatedHits.load(\"items\");\nawait context.sync();\n\nconst updatedRange = updatedHits.items[0];\nconst tailHits = updatedRange.search(\", 2015\", { matchCase: true });\ntailHits.load(\"items\");\nawait context.sync();\n\nconst insertionPoint = tailHits.items[0].getRange(\"Start\");\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document's auto-managed \"_GoBack\" bookmark currently wraps the\n# class-diagram picture (left over from whatever was typed/pasted there\n# previously). We are about to make the \"real\" last edit -- changing the\n# date from \"February 4, 2015\" to \"February 25, 2015\" -- so the bookmark\n# needs to move to mark that new edit position instead.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the stale \"_GoBack\" bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Find the date text, then find the day \"4\" within it and replace with \"25\".\n$dateRange = $d.Content\n$null = $dateRange.Find.Execute(\"February 4, 2015\")\n\n$dayRange = $dateRange.Duplicate\n$null = $dayRange.Find.Execute(\"4\", $false, $false, $false, $false, $false, $true, 1, $false, \"25\", 2)\n\n# 3) Collapse the (now \"25\") range to its end point -- i.e. immediately\n#    before the \", 2015\" that follows -- and drop a collapsed \"_GoBack\"\n#    bookmark right there, marking the new edit position.\n$dayRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $dayRange)\n"}
